# Auto-generated update of market-price / profit columns (H-N) across all 8 job sheets.
# Mirrors a scheduled data refresh: most cells get new values; a few cells are
# cleared (diff removes the <c> element) and a few new cells are added where a
# profit column did not previously have a value.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1978
$ws.Range("I2").Value = 795
$ws.Range("K2").Value = 795
$ws.Range("M2").Value = -682
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H76").Value = 4480.6113
$ws.Range("I76").Value = 3788.923
$ws.Range("K76").Value = 3788.923
$ws.Range("M76").Value = -3473.923
$ws.Range("H79").Value = 4480.6113
$ws.Range("I79").Value = 3788.923
$ws.Range("K79").Value = 3788.923
$ws.Range("M79").Value = -2696.923
$ws.Range("H80").Value = 17841.584
$ws.Range("I80").Value = 12662
$ws.Range("J80").Value = 20431.375
$ws.Range("K80").Value = 37986
$ws.Range("L80").Value = 61294.125
$ws.Range("M80").Value = -36988
$ws.Range("N80").Value = -63290.125
$ws.Range("H83").Value = 17841.584
$ws.Range("I83").Value = 12662
$ws.Range("J83").Value = 20431.375
$ws.Range("K83").Value = 113958
$ws.Range("L83").Value = 183882.375
$ws.Range("M83").Value = -108966
$ws.Range("N83").Value = -193866.375
$ws.Range("H86").Value = 2593
$ws.Range("I86").Value = 2413.6667
$ws.Range("K86").Value = 2413.6667
$ws.Range("M86").Value = -1290.6667
$ws.Range("H89").Value = 2593
$ws.Range("I89").Value = 2413.6667
$ws.Range("K89").Value = 12068.3335
$ws.Range("M89").Value = -6452.333500000001
$ws.Range("H98").Value = 3120.875
$ws.Range("I98").Value = 2783.28
$ws.Range("K98").Value = 2783.28
$ws.Range("M98").Value = -1285.28
$ws.Range("H119").Value = 1909.2858
$ws.Range("J119").Value = 1909.2858
$ws.Range("L119").Value = 5727.857400000001
$ws.Range("N119").Value = -15403.8574
$ws.Range("H122").Value = 3120.875
$ws.Range("I122").Value = 2783.28
$ws.Range("K122").Value = 8349.84
$ws.Range("M122").Value = -5899.84
$ws.Range("H132").Value = 1915730.4
$ws.Range("I132").Value = 2176490.8
$ws.Range("J132").Value = 3488
$ws.Range("K132").Value = 6529472.399999999
$ws.Range("L132").Value = 10464
$ws.Range("M132").Value = -6526942.399999999
$ws.Range("N132").Value = -15524
$ws.Range("H135").Value = 20283.355
$ws.Range("I135").Value = 1447.9166
$ws.Range("K135").Value = 13031.2494
$ws.Range("M135").Value = -10496.2494
$ws.Range("H137").Value = 28826.834
$ws.Range("I137").Value = 32192.3
$ws.Range("K137").Value = 96576.89999999999
$ws.Range("M137").Value = -94026.89999999999
$ws.Range("H138").Value = 1805.8286
$ws.Range("I138").Value = 1455.4138
$ws.Range("J138").Value = 3499.5
$ws.Range("K138").Value = 4366.2414
$ws.Range("L138").Value = 10498.5
$ws.Range("M138").Value = 773.7586000000001
$ws.Range("N138").Value = -20778.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2351.5
$ws.Range("I2").Value = 2124
$ws.Range("K2").Value = 2124
$ws.Range("M2").Value = -2011
$ws.Range("H32").Value = 36662.812
$ws.Range("I32").Value = 31986.553
$ws.Range("J32").Value = 54432.6
$ws.Range("K32").Value = 31986.553
$ws.Range("L32").Value = 54432.6
$ws.Range("M32").Value = -31699.553
$ws.Range("N32").Value = -55006.6
$ws.Range("H45").Value = 9193.875
$ws.Range("I45").Value = 8341
$ws.Range("K45").Value = 8341
$ws.Range("M45").Value = -7964
$ws.Range("H61").Value = 2517.1572
$ws.Range("I61").Value = 717.98303
$ws.Range("J61").Value = 12167.272
$ws.Range("K61").Value = 717.98303
$ws.Range("L61").Value = 12167.272
$ws.Range("M61").Value = -505.98303
$ws.Range("N61").Value = -12591.272
$ws.Range("H74").Value = 191913.31
$ws.Range("I74").Value = 215543.83
$ws.Range("K74").Value = 215543.83
$ws.Range("M74").Value = -214669.83
$ws.Range("H77").Value = 191913.31
$ws.Range("I77").Value = 215543.83
$ws.Range("K77").Value = 1077719.15
$ws.Range("M77").Value = -1073351.15
$ws.Range("H116").Value = 2351.5
$ws.Range("I116").Value = 2124
$ws.Range("K116").Value = 2124
$ws.Range("M116").Value = 170
$ws.Range("H132").Value = 1423.4395
$ws.Range("I132").Value = 1190.6415
$ws.Range("J132").Value = 2372.5386
$ws.Range("K132").Value = 3571.9245
$ws.Range("L132").Value = 7117.6158
$ws.Range("M132").Value = -1041.9245
$ws.Range("N132").Value = -12177.6158
$ws.Range("H136").Value = 2517.1572
$ws.Range("I136").Value = 717.98303
$ws.Range("J136").Value = 12167.272
$ws.Range("K136").Value = 2153.94909
$ws.Range("L136").Value = 36501.81600000001
$ws.Range("M136").Value = 396.0509099999999
$ws.Range("N136").Value = -41601.81600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2351.5
$ws.Range("I3").Value = 2124
$ws.Range("K3").Value = 2124
$ws.Range("M3").Value = -2010
$ws.Range("H59").Value = 70354.5
$ws.Range("J59").Value = 100000
$ws.Range("L59").Value = 100000
$ws.Range("N59").Value = -101694
$ws.Range("H134").Value = 3372.4893
$ws.Range("I134").Value = 1836.3793
$ws.Range("J134").Value = 5847.3335
$ws.Range("K134").Value = 5509.1379
$ws.Range("L134").Value = 17542.0005
$ws.Range("M134").Value = -2974.1379
$ws.Range("N134").Value = -22612.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 5499.875
$ws.Range("I19").Value = 1274.75
$ws.Range("K19").Value = 1274.75
$ws.Range("M19").Value = -1104.75
$ws.Range("H22").Value = 1904.4
$ws.Range("J22").Value = 2262.25
$ws.Range("L22").Value = 2262.25
$ws.Range("N22").Value = -2962.25
$ws.Range("H24").Value = 5499.875
$ws.Range("I24").Value = 1274.75
$ws.Range("K24").Value = 1274.75
$ws.Range("M24").Value = -1104.75
$ws.Range("H31").Value = 5559202.5
$ws.Range("I31").Value = 6670043
$ws.Range("J31").Value = 4999.3335
$ws.Range("K31").Value = 6670043
$ws.Range("L31").Value = 4999.3335
$ws.Range("M31").Value = -6669748
$ws.Range("N31").Value = -5589.3335
$ws.Range("H34").Value = 5559202.5
$ws.Range("I34").Value = 6670043
$ws.Range("J34").Value = 4999.3335
$ws.Range("K34").Value = 6670043
$ws.Range("L34").Value = 4999.3335
$ws.Range("M34").Value = -6669841
$ws.Range("N34").Value = -5403.3335
$ws.Range("H51").Value = 37000
$ws.Range("J51").Value = 37000
$ws.Range("L51").Value = 37000
$ws.Range("N51").Value = -38472
$ws.Range("H58").Value = 2078.3333
$ws.Range("I58").Value = 1878
$ws.Range("J58").Value = 3080
$ws.Range("K58").Value = 1878
$ws.Range("L58").Value = 3080
$ws.Range("M58").Value = -1675
$ws.Range("N58").Value = -3486
$ws.Range("H61").Value = 37000
$ws.Range("J61").Value = 37000
$ws.Range("L61").Value = 37000
$ws.Range("N61").Value = -37696
$ws.Range("H64").Value = 35000
$ws.Range("J64").Value = 35000
$ws.Range("L64").Value = 35000
$ws.Range("N64").Value = -35496
$ws.Range("H67").Value = 35000
$ws.Range("J67").Value = 35000
$ws.Range("L67").Value = 35000
$ws.Range("N67").Value = -36716
$ws.Range("H93").Value = 31851.5
$ws.Range("I93").Value = 31851.5
$ws.Range("K93").Value = 31851.5
$ws.Range("M93").Value = -29979.5
$ws.Range("H94").Value = 641.125
$ws.Range("I94").Value = 797.25
$ws.Range("J94").Value = 485
$ws.Range("K94").Value = 797.25
$ws.Range("L94").Value = 485
$ws.Range("M94").Value = -346.25
$ws.Range("N94").Value = -1387
$ws.Range("H132").Value = 42663.875
$ws.Range("I132").Value = 46362.863
$ws.Range("J132").Value = 1975
$ws.Range("K132").Value = 139088.589
$ws.Range("L132").Value = 5925
$ws.Range("M132").Value = -136558.589
$ws.Range("N132").Value = -10985
$ws.Range("H134").Value = 1593.7142
$ws.Range("I134").Value = 1233.091
$ws.Range("J134").Value = 2916
$ws.Range("K134").Value = 3699.273
$ws.Range("L134").Value = 8748
$ws.Range("M134").Value = -1164.273
$ws.Range("N134").Value = -13818
$ws.Range("H136").Value = 2078.3333
$ws.Range("I136").Value = 1878
$ws.Range("J136").Value = 3080
$ws.Range("K136").Value = 5634
$ws.Range("L136").Value = 9240
$ws.Range("M136").Value = -3084
$ws.Range("N136").Value = -14340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 56193784
$ws.Range("I4").Value = 56686812
$ws.Range("J4").Value = 55723172
$ws.Range("K4").Value = 170060436
$ws.Range("L4").Value = 167169516
$ws.Range("M4").Value = -170060324
$ws.Range("N4").Value = -167169740
$ws.Range("H5").Value = 1601.1111
$ws.Range("I5").Value = 1306.2307
$ws.Range("K5").Value = 3918.6921
$ws.Range("M5").Value = -3806.6921
$ws.Range("H11").Value = 1617781.2
$ws.Range("I11").Value = 2101017.8
$ws.Range("K11").Value = 6303053.399999999
$ws.Range("M11").Value = -6302913.399999999
$ws.Range("H22").Value = 4525.615
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 4525.615
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 13576.845
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -13914.845
$ws.Range("H27").Value = 4525.615
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 4525.615
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 13576.845
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -13780.845
$ws.Range("H75").Value = 3626.7742
$ws.Range("I75").Value = 1573
$ws.Range("J75").Value = 4021.7307
$ws.Range("K75").Value = 4719
$ws.Range("L75").Value = 12065.1921
$ws.Range("M75").Value = -3721
$ws.Range("N75").Value = -14061.1921
$ws.Range("H78").Value = 3626.7742
$ws.Range("I78").Value = 1573
$ws.Range("J78").Value = 4021.7307
$ws.Range("K78").Value = 14157
$ws.Range("L78").Value = 36195.5763
$ws.Range("M78").Value = -9165
$ws.Range("N78").Value = -46179.5763
$ws.Range("H92").Value = 1837.8889
$ws.Range("I92").Value = 2616.1667
$ws.Range("J92").Value = 281.33334
$ws.Range("K92").Value = 7848.500100000001
$ws.Range("L92").Value = 844.0000200000001
$ws.Range("M92").Value = -6600.500100000001
$ws.Range("N92").Value = -3340.00002
$ws.Range("H117").Value = 510.83334
$ws.Range("I117").Value = 595
$ws.Range("J117").Value = 90
$ws.Range("K117").Value = 1785
$ws.Range("L117").Value = 270
$ws.Range("M117").Value = 1657
$ws.Range("N117").Value = -7154
$ws.Range("H122").Value = 700.57574
$ws.Range("I122").Value = 1074.4
$ws.Range("J122").Value = 633.8214
$ws.Range("K122").Value = 9669.6
$ws.Range("L122").Value = 5704.3926
$ws.Range("M122").Value = -7219.6
$ws.Range("N122").Value = -10604.3926
$ws.Range("H134").Value = 659
$ws.Range("I134").Value = 659
$ws.Range("K134").Value = 1977
$ws.Range("M134").Value = 3093
$ws.Range("H135").Value = 1601.1111
$ws.Range("I135").Value = 1306.2307
$ws.Range("K135").Value = 11756.0763
$ws.Range("M135").Value = -9221.076300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1045.8572
$ws.Range("J2").Value = 712
$ws.Range("L2").Value = 712
$ws.Range("N2").Value = -938
$ws.Range("H80").Value = 8551.516
$ws.Range("I80").Value = 5570.75
$ws.Range("K80").Value = 5570.75
$ws.Range("M80").Value = -4572.75
$ws.Range("H83").Value = 8551.516
$ws.Range("I83").Value = 5570.75
$ws.Range("K83").Value = 27853.75
$ws.Range("M83").Value = -22861.75
$ws.Range("H102").Value = 40858.152
$ws.Range("I102").Value = 47666.184
$ws.Range("J102").Value = 3414
$ws.Range("K102").Value = 47666.184
$ws.Range("L102").Value = 3414
$ws.Range("M102").Value = -46044.184
$ws.Range("N102").Value = -6658
$ws.Range("H126").Value = 4045.7778
$ws.Range("I126").Value = 2066.5
$ws.Range("K126").Value = 6199.5
$ws.Range("M126").Value = -3729.5
$ws.Range("H132").Value = 1395.1285
$ws.Range("I132").Value = 1372.356
$ws.Range("J132").Value = 1517.2727
$ws.Range("K132").Value = 4117.068
$ws.Range("L132").Value = 4551.8181
$ws.Range("M132").Value = -1587.068
$ws.Range("N132").Value = -9611.8181

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2663.2144
$ws.Range("I40").Value = 2654.5652
$ws.Range("K40").Value = 2654.5652
$ws.Range("M40").Value = -2518.5652
$ws.Range("H46").Value = 4475.857
$ws.Range("I46").Value = 2114.2856
$ws.Range("J46").Value = 6837.4287
$ws.Range("K46").Value = 2114.2856
$ws.Range("L46").Value = 6837.4287
$ws.Range("M46").Value = -1926.2856
$ws.Range("N46").Value = -7213.4287
$ws.Range("H61").Value = 1135.3636
$ws.Range("I61").Value = 1069.8823
$ws.Range("K61").Value = 1069.8823
$ws.Range("M61").Value = -867.8823
$ws.Range("H108").Value = 72899
$ws.Range("J108").Value = 72899
$ws.Range("L108").Value = 72899
$ws.Range("N108").Value = -80579
$ws.Range("H113").Value = 1135.3636
$ws.Range("I113").Value = 1069.8823
$ws.Range("K113").Value = 1069.8823
$ws.Range("M113").Value = 1100.1177
$ws.Range("H132").Value = 2847.4827
$ws.Range("I132").Value = 2847.4827
$ws.Range("K132").Value = 8542.4481
$ws.Range("M132").Value = -6012.4481
$ws.Range("H136").Value = 2177.3542
$ws.Range("I136").Value = 1951.3695
$ws.Range("K136").Value = 5854.1085
$ws.Range("M136").Value = -3304.1085

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 89999.336
$ws.Range("I75").Value = 69999
$ws.Range("K75").Value = 69999
$ws.Range("M75").Value = -69063
$ws.Range("H78").Value = 89999.336
$ws.Range("I78").Value = 69999
$ws.Range("K78").Value = 209997
$ws.Range("M78").Value = -205317
$ws.Range("H132").Value = 28534.275
$ws.Range("I132").Value = 37613.617
$ws.Range("K132").Value = 112840.851
$ws.Range("M132").Value = -110310.851
$ws.Range("H136").Value = 9547.383
$ws.Range("I136").Value = 10655.879
$ws.Range("J136").Value = 4670
$ws.Range("K136").Value = 31967.637
$ws.Range("L136").Value = 14010
$ws.Range("M136").Value = -29417.637
$ws.Range("N136").Value = -19110
